$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.584.83'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.69'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.90'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4604'
$ws.Range("E7").Value = '  -1.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3860'
$ws.Range("E8").Value = '  -2.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.58'
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07868'
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +2.19%  '

$ws.Range("E12").Value = '  -2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.885.68'
$ws.Range("E13").Value = '  -1.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.056'
$ws.Range("E14").Value = '  +0.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.698'
$ws.Range("E15").Value = '  -0.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06975'
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.43'
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.009'
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001003'
$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.19'
$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.008'
$ws.Range("E21").Value = '  +0.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.582.99'
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.331'
$ws.Range("E23").Value = '  -0.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.96'
$ws.Range("E24").Value = '  -1.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.101.11'
$ws.Range("E25").Value = '  -1.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.057'
$ws.Range("E26").Value = '  -3.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.55'
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.47'
$ws.Range("E28").Value = '  +0.14%  '

$ws.Range("E29").Value = '  +1.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.955'
$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.42'
$ws.Range("E31").Value = '  -1.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09340'
$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("E33").Value = '  -2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.316'
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.336'
$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("E36").Value = '  -2.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05768'
$ws.Range("E37").Value = '  -2.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.164'
$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.983'
$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("E40").Value = '  -2.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5679'
$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1794'
$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.728'
$ws.Range("E43").Value = '  -2.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.83'
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5364'
$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07138'
$ws.Range("E46").Value = '  -1.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.146'
$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.840'
$ws.Range("E48").Value = '  -0.64%  '

$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.118'
$ws.Range("E49").Value = '  -2.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.24'
$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.489'
$ws.Range("E51").Value = '  +5.05%  '
